$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - leading apostrophe forces these numeric-looking
# strings to stay stored as text, matching the source data's inline-string
# format instead of being auto-coerced into numbers.
$ws.Range("D2").Value  = "'246.76"
$ws.Range("D3").Value  = "'26.28"
$ws.Range("D4").Value  = "'5.065"
$ws.Range("D6").Value  = "'6.498"
$ws.Range("D7").Value  = "'3.045"
$ws.Range("D9").Value  = "'0.8409"
$ws.Range("D10").Value = "'0.1345"
$ws.Range("D11").Value = "'0.02810"
$ws.Range("D12").Value = "'0.09380"
$ws.Range("D13").Value = "'0.001519"
$ws.Range("D14").Value = "'0.0005967"
$ws.Range("D15").Value = "'0.006185"
$ws.Range("D16").Value = "'3.552"
$ws.Range("D19").Value = "'0.06951"
$ws.Range("D20").Value = "'0.03117"
$ws.Range("D22").Value = "'3.740"
$ws.Range("D23").Value = "'0.04671"
$ws.Range("D24").Value = "'0.1374"
$ws.Range("D25").Value = "'0.001246"
$ws.Range("D26").Value = "'0.004614"
$ws.Range("D28").Value = "'0.0001389"
$ws.Range("D41").Value = "'0.006177"
$ws.Range("D42").Value = "'0.1056"
$ws.Range("D44").Value = "'0.008988"
$ws.Range("D45").Value = "'0.00005287"
$ws.Range("D47").Value = "'0.1598"

# Volume(1h) / "Best-Worst in 24h" label (column E) updates
$ws.Range("E14").Value = "13OneONE"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
